$d = $word.ActiveDocument

# --- 1. Tidy up the "Agenda item #3" discussion paragraph -----------------
# In the original document this paragraph is split across four runs with
# w:proofErr gramStart/gramEnd markers bracketing the word "Also" (a
# leftover from Word's grammar checker). The edit collapses that into two
# runs (no proofErr markers), keeping the final "employees can have a
# maximum..." run separate.
$discussionPara = $d.Paragraphs.Item(17)
$fullRange = $discussionPara.Range
# Exclude the trailing paragraph mark so the <w:p> itself (and its
# paraId/rsid attributes) is left untouched; only the runs inside change.
$bodyRange = $d.Range($fullRange.Start, $fullRange.End - 1)

$bodyXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body><w:p>' + `
            '<w:r><w:t xml:space="preserve"> The customer wants the administration to be able to access and edit pretty much everything, the manager should only be able to look at statistics and the working shift schedule, and the workers should only be able to look at their own schedules. Also the </w:t></w:r>' + `
            '<w:r><w:t>employees can have a maximum of 2 shifts per day and a shift can have a maximum of 10 employees.</w:t></w:r>' + `
          '</w:p></w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$bodyRange.InsertXML($bodyXml)

# --- 2. Add the Location/Participants header block -------------------------
# The document originally starts with:
#   P1: "Date/Time of Meeting: 14/02/2022"   (Heading1)
#   P2: (empty paragraph, no explicit style)
#   P3: "Agenda item #1: ..."
#   ...
#
# The edit restructures the minutes by turning the single empty
# paragraph (P2) into a proper header block:
#   - "Location: R10 building, Fontys"                         (Heading2)
#   - "Participants: Rositsa Nikolova, Daniil Blagoev, ..."     (Heading2, nl-NL)
#   - an empty paragraph carrying the nl-NL paragraph-mark language
#
# Replace that paragraph's range with the three new paragraphs via a
# raw OOXML injection so exact formatting (styles + nl-NL language tag)
# is reproduced precisely.
$locationParagraph = $d.Paragraphs.Item(2)
$target = $locationParagraph.Range

$headerXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Location: R10 building, Fontys</w:t></w:r></w:p>' + `
            '<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t>Participants: Rositsa Nikolova, Daniil Blagoev, Jakub Jelinek, Rens van den Elzen</w:t></w:r></w:p>' + `
            '<w:p><w:pPr><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr></w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$target.InsertXML($headerXml)
